# Auto-generated edit script: updates market price / profit values
# in the per-job leve profit tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 449.75
$ws.Range("I9").Value = 319.6
$ws.Range("K9").Value = 319.6
$ws.Range("M9").Value = -150.6
$ws.Range("H12").Value = 462.5625
$ws.Range("I12").Value = 486.73334
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 486.73334
$ws.Range("L12").Value = 100
$ws.Range("M12").Value = -316.73334
$ws.Range("N12").Value = -440
$ws.Range("H40").Value = 3567.3076
$ws.Range("I40").Value = 2641.625
$ws.Range("J40").Value = 5048.4
$ws.Range("K40").Value = 2641.625
$ws.Range("L40").Value = 5048.4
$ws.Range("M40").Value = -2466.625
$ws.Range("N40").Value = -5398.4
$ws.Range("H80").Value = 28256.75
$ws.Range("J80").Value = 28256.75
$ws.Range("L80").Value = 84770.25
$ws.Range("N80").Value = -86766.25
$ws.Range("H83").Value = 28256.75
$ws.Range("J83").Value = 28256.75
$ws.Range("L83").Value = 254310.75
$ws.Range("N83").Value = -264294.75
$ws.Range("I132").Value = 6838263.5
$ws.Range("K132").Value = 20514790.5
$ws.Range("M132").Value = -20512260.5
$ws.Range("H135").Value = 16478.363
$ws.Range("I135").Value = 773.52
$ws.Range("K135").Value = 6961.68
$ws.Range("M135").Value = -4426.68
$ws.Range("H137").Value = 13860.759
$ws.Range("I137").Value = 17509.45
$ws.Range("J137").Value = 5752.5557
$ws.Range("K137").Value = 52528.35000000001
$ws.Range("L137").Value = 17257.6671
$ws.Range("M137").Value = -49978.35000000001
$ws.Range("N137").Value = -22357.6671
$ws.Range("H138").Value = 3718
$ws.Range("I138").Value = 3224.158
$ws.Range("K138").Value = 9672.474
$ws.Range("M138").Value = -4532.474
$ws.Range("H141").Value = 2323.3333
$ws.Range("I141").Value = 2323.3333
$ws.Range("K141").Value = 6969.999899999999
$ws.Range("M141").Value = -1789.999899999999

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20895.732
$ws.Range("I32").Value = 22971.28
$ws.Range("K32").Value = 22971.28
$ws.Range("M32").Value = -22684.28
$ws.Range("H61").Value = 5205.5806
$ws.Range("I61").Value = 1632.7894
$ws.Range("K61").Value = 1632.7894
$ws.Range("M61").Value = -1420.7894
$ws.Range("H63").Value = 3882.9333
$ws.Range("I63").Value = 1965
$ws.Range("K63").Value = 1965
$ws.Range("M63").Value = -1279
$ws.Range("H66").Value = 3882.9333
$ws.Range("I66").Value = 1965
$ws.Range("K66").Value = 9825
$ws.Range("M66").Value = -6393
$ws.Range("H132").Value = 1411.44
$ws.Range("I132").Value = 1182.6111
$ws.Range("K132").Value = 3547.8333
$ws.Range("M132").Value = -1017.8333
$ws.Range("H136").Value = 5205.5806
$ws.Range("I136").Value = 1632.7894
$ws.Range("K136").Value = 4898.3682
$ws.Range("M136").Value = -2348.3682

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 746.36
$ws.Range("J80").Value = 667.7143
$ws.Range("L80").Value = 667.7143
$ws.Range("N80").Value = -2663.7143
$ws.Range("H83").Value = 746.36
$ws.Range("J83").Value = 667.7143
$ws.Range("L83").Value = 3338.5715
$ws.Range("N83").Value = -13322.5715
$ws.Range("H94").Value = 5295.522
$ws.Range("I94").Value = 6223.0586
$ws.Range("K94").Value = 6223.0586
$ws.Range("M94").Value = -5772.0586
$ws.Range("H134").Value = 2975.4167
$ws.Range("I134").Value = 2235.0527
$ws.Range("K134").Value = 6705.158100000001
$ws.Range("M134").Value = -4170.158100000001

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 849.5
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 699
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 699
$ws.Range("M2").Value = -887
$ws.Range("N2").Value = -925
$ws.Range("H3").Value = 2000
$ws.Range("J3").Value = 2000
$ws.Range("L3").Value = 2000
$ws.Range("N3").Value = -2226
$ws.Range("H6").Value = 1588124
$ws.Range("I6").Value = 5498.625
$ws.Range("K6").Value = 5498.625
$ws.Range("M6").Value = -5385.625
$ws.Range("H16").Value = 1855.7273
$ws.Range("I16").Value = 1523.0769
$ws.Range("K16").Value = 1523.0769
$ws.Range("M16").Value = -1236.0769
$ws.Range("H31").Value = 10004970
$ws.Range("I31").Value = 25001500
$ws.Range("J31").Value = 7283
$ws.Range("K31").Value = 25001500
$ws.Range("L31").Value = 7283
$ws.Range("M31").Value = -25001205
$ws.Range("N31").Value = -7873
$ws.Range("H34").Value = 10004970
$ws.Range("I34").Value = 25001500
$ws.Range("J34").Value = 7283
$ws.Range("K34").Value = 25001500
$ws.Range("L34").Value = 7283
$ws.Range("M34").Value = -25001298
$ws.Range("N34").Value = -7687
$ws.Range("H52").Value = 95000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 95000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 95000
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -95588
$ws.Range("H113").Value = 1855.7273
$ws.Range("I113").Value = 1523.0769
$ws.Range("K113").Value = 1523.0769
$ws.Range("M113").Value = 646.9231
$ws.Range("H133").Value = 75318.5
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 75318.5
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 75318.5
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -80378.5
$ws.Range("H134").Value = 1343.5938
$ws.Range("I134").Value = 1343.5938
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4030.7814
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1495.7814
$ws.Range("N134").ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1500744.9
$ws.Range("I11").Value = 1750762.2
$ws.Range("J11").Value = 640.5
$ws.Range("K11").Value = 5252286.6
$ws.Range("L11").Value = 1921.5
$ws.Range("M11").Value = -5252146.6
$ws.Range("N11").Value = -2201.5
$ws.Range("H60").Value = 2803.9285
$ws.Range("J60").Value = 4045.111
$ws.Range("L60").Value = 12135.333
$ws.Range("N60").Value = -12637.333
$ws.Range("H81").Value = 3963.3333
$ws.Range("J81").Value = 10000
$ws.Range("L81").Value = 30000
$ws.Range("N81").Value = -32246
$ws.Range("H84").Value = 3963.3333
$ws.Range("J84").Value = 10000
$ws.Range("L84").Value = 90000
$ws.Range("N84").Value = -101232
$ws.Range("H105").Value = 7233.3335
$ws.Range("I105").Value = 7200
$ws.Range("K105").Value = 21600
$ws.Range("M105").Value = -18979
$ws.Range("H107").Value = 613.75
$ws.Range("I107").Value = 406.5
$ws.Range("K107").Value = 1219.5
$ws.Range("M107").Value = 700.5

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 12000
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H62").Value = 60000
$ws.Range("J62").Value = 60000
$ws.Range("L62").Value = 60000
$ws.Range("N62").Value = -61372
$ws.Range("H65").Value = 60000
$ws.Range("J65").Value = 60000
$ws.Range("L65").Value = 180000
$ws.Range("N65").Value = -186864
$ws.Range("H107").Value = 333.55554
$ws.Range("J107").Value = 412
$ws.Range("L107").Value = 412
$ws.Range("N107").Value = -4252
$ws.Range("H132").Value = 2437.5715
$ws.Range("I132").Value = 2450.8542
$ws.Range("J132").Value = 1800
$ws.Range("K132").Value = 7352.562600000001
$ws.Range("L132").Value = 5400
$ws.Range("M132").Value = -4822.562600000001
$ws.Range("N132").Value = -10460
$ws.Range("H136").Value = 6771
$ws.Range("J136").Value = 6771
$ws.Range("L136").Value = 20313
$ws.Range("N136").Value = -25413
$ws.Range("H137").Value = 80707.164
$ws.Range("J137").Value = 80707
$ws.Range("L137").Value = 80707
$ws.Range("N137").Value = -90907

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1243.8096
$ws.Range("I16").Value = 1060.0714
$ws.Range("J16").Value = 1611.2858
$ws.Range("K16").Value = 1060.0714
$ws.Range("L16").Value = 1611.2858
$ws.Range("M16").Value = -890.0714
$ws.Range("N16").Value = -1951.2858
$ws.Range("H34").Value = 79999.5
$ws.Range("I34").Value = 79999
$ws.Range("K34").Value = 79999
$ws.Range("M34").Value = -79827
$ws.Range("H132").Value = 2373.2856
$ws.Range("I132").Value = 1788.3928
$ws.Range("J132").Value = 4712.857
$ws.Range("K132").Value = 5365.178400000001
$ws.Range("L132").Value = 14138.571
$ws.Range("M132").Value = -2835.178400000001
$ws.Range("N132").Value = -19198.571

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H132").Value = 1311.2821
$ws.Range("I132").Value = 768.41174
$ws.Range("J132").Value = 5002.8
$ws.Range("K132").Value = 2305.23522
$ws.Range("L132").Value = 15008.4
$ws.Range("M132").Value = 224.76478
$ws.Range("N132").Value = -20068.4

